$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 526, shifting existing rows 526-573 down to 527-574.
$ws.Rows("526:526").Insert()

# Populate the newly inserted row 526 with the new record.
$ws.Range("A526").Value = 5
$ws.Range("B526").Value = "Macroferia Regional de Talca"
$ws.Range("C526").Value = "Maule"
$ws.Range("D526").Value = 45106
$ws.Range("E526").Value = 7
$ws.Range("F526").Value = 100112023
$ws.Range("G526").Value = "Brócoli"
$ws.Range("H526").Value = "Sin especificar"
$ws.Range("I526").Value = "Primera"
$ws.Range("J526").Value = 5000
$ws.Range("K526").Value = 500
$ws.Range("L526").Value = 500
$ws.Range("M526").Value = 500
$ws.Range("N526").Value = "$/unidad"
$ws.Range("O526").Value = "Región del Maule"
$ws.Range("P526").Value = 500
$ws.Range("Q526").Value = 1
$ws.Range("R526").Value = "Hortaliza"
